# Updated symbol list on Fri Jan 13 08:42:27 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Rows 2-6: price / volume refresh only
Set-TextCell "D2" "287.49"
Set-TextCell "E2" "1.45%"

Set-TextCell "D3" "29.59"
Set-TextCell "E3" "3.38%"

Set-TextCell "D4" "5.121"
Set-TextCell "E4" "1.04%"

Set-TextCell "D5" "0.06692"
Set-TextCell "E5" "3.26%"

Set-TextCell "D6" "7.330"
Set-TextCell "E6" "1.64%"

# Rows 7-19: coin list shifted down by one (new entrant at row 7),
# plus refreshed price / volume figures
Set-TextCell "B7" "GateToken"
Set-TextCell "C7" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell "D7" "3.405"
Set-TextCell "E7" "1.12%"

Set-TextCell "B8" "FTXToken"
Set-TextCell "C8" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell "D8" "1.362"
Set-TextCell "E8" "-1.15%"

Set-TextCell "B9" "MXToken"
Set-TextCell "C9" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell "D9" "0.9128"
Set-TextCell "E9" "0.20%"

Set-TextCell "B10" "WazirX"
Set-TextCell "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell "D10" "0.1588"
Set-TextCell "E10" "2.48%"

Set-TextCell "B11" "LiechtensteinCryptoassetsExchange"
Set-TextCell "C11" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell "D11" "0.06751"
Set-TextCell "E11" "3.62%"

Set-TextCell "B12" "MandalaExchangeToken"
Set-TextCell "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell "D12" "0.07729"
Set-TextCell "E12" "1.59%"

Set-TextCell "B13" "BitrueCoin"
Set-TextCell "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell "D13" "0.02933"
Set-TextCell "E13" "6.56%"

Set-TextCell "B14" "BitMartToken"
Set-TextCell "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell "D14" "0.08976"
Set-TextCell "E14" "0.29%"

Set-TextCell "B15" "BitForexToken"
Set-TextCell "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell "D15" "0.001579"
Set-TextCell "E15" "-0.76%"

Set-TextCell "B16" "CoinExToken"
Set-TextCell "C16" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextCell "D16" "0.04491"
Set-TextCell "E16" "0.86%"

Set-TextCell "B17" "One"
Set-TextCell "C17" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell "D17" "0.0006461"
Set-TextCell "E17" "1.85%"

Set-TextCell "B18" "TigerCash"
Set-TextCell "C18" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell "D18" "0.006273"
Set-TextCell "E18" "3.01%"

Set-TextCell "B19" "LEO"
Set-TextCell "C19" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell "D19" "3.441"
Set-TextCell "E19" "-0.26%"

# Rows 20-28: price / volume refresh only
Set-TextCell "D20" "2.221"
Set-TextCell "E20" "-0.90%"

Set-TextCell "E21" "0.77%"

Set-TextCell "E22" "-2.37%"

Set-TextCell "D23" "4.069"
Set-TextCell "E23" "1.73%"

Set-TextCell "E25" "0.72%"

Set-TextCell "D26" "0.004113"
Set-TextCell "E26" "-5.61%"

Set-TextCell "E27" "-0.14%"

Set-TextCell "D28" "0.0001617"
Set-TextCell "E28" "-1.16%"

# Rows 40-47: price / volume refresh only
Set-TextCell "D40" "0.04263"
Set-TextCell "E40" "3.58%"

Set-TextCell "D41" "0.006763"
Set-TextCell "E41" "1.78%"

Set-TextCell "D42" "0.1239"
Set-TextCell "E42" "0.70%"

Set-TextCell "E43" "7.57%"

Set-TextCell "D44" "0.01329"
Set-TextCell "E44" "6.64%"

Set-TextCell "D45" "0.00005696"
Set-TextCell "E45" "5.47%"

Set-TextCell "E46" "1.81%"

Set-TextCell "D47" "0.01306"
Set-TextCell "E47" "-29.42%"
